# Automatic update of files.
# Update "Förändrad" (column C) date for all data rows (+1 day),
# and re-sort rows 3-9 ascending by "Datum" (column B), carrying along
# the related "Beteckning" (A) and "Area (ha)" (G) values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 keeps its place; only the "Förändrad" date advances by one day.
$ws.Range("C2").Value = 46077

# New ordering (rows 3-9), sorted ascending by column B (date), taken
# from the current contents of rows 3-9.
$data = @(
    @{ A = "A 54782-2022"; B = 44883; C = 46077; G = 5.5 },
    @{ A = "A 34310-2024"; B = 45524; C = 46077; G = 4.8 },
    @{ A = "A 843-2024";   B = 45300; C = 46077; G = 0.8 },
    @{ A = "A 844-2024";   B = 45300; C = 46077; G = 1.2 },
    @{ A = "A 45983-2023"; B = 45196; C = 46077; G = 0.6 },
    @{ A = "A 17908-2021"; B = 44301; C = 46077; G = 0.9 },
    @{ A = "A 25617-2024"; B = 45463; C = 46077; G = 2.3 }
)

$r = 3
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    $ws.Cells.Item($r, 7).Value = $row.G
    $r = $r + 1
}
